# Bug fix: the English string "Initial & Final Surveillance Diagnosis" (row 74)
# doesn't translate correctly because of the "&" character. A new row is added
# right after it with an "&"-free variant ("Initial and Final Surveillance
# Diagnosis") that carries the existing Vietnamese translation, while the
# original row's translation is reset to "TBT" (to be translated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the Vietnamese translation currently sitting in B74; it belongs to
# the new "and" row we are about to insert.
$existingTranslation = $ws.Range("B74").Value2

# Insert a new blank row at position 75; everything from the old row 75 down
# shifts down by one row (matches the dimension growing from B181 to B182).
$ws.Rows("75:75").Insert()

# Populate the newly inserted row 75.
$ws.Range("A75").Value2 = "Initial and Final Surveillance Diagnosis"
$ws.Range("B75").Value2 = $existingTranslation

# Reset the translation on the original row 74 (the "&" variant).
$ws.Range("B74").Value2 = "TBT"
